$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = "missing_filepaths"
$ws.Range("E11").Value = "AAGTGCAGCTCGTCCGGCGT"

$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("E11").Select()
